$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 39 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A39").Value = "Vragen over samenwerking"
$logs.Range("B39").Value = "mailmind.test@zohomail.eu"
$logs.Range("C39").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D39").Value = "Overig"
$logs.Range("F39").Value = "2025-06-17 21:59:06"
$logs.Range("G39").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges to include the new row ---
$dFc = $logs.Range("D2:D38").FormatConditions.Item(1)
$dFc.ModifyAppliesToRange($logs.Range("D2:D39"))

$gFc = $logs.Range("G2:G38").FormatConditions.Item(1)
$gFc.ModifyAppliesToRange($logs.Range("G2:G39"))

# --- Dashboard sheet: update "Overig" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 11
